# CIV-14127 Updated typo in send to other party template
#
# The template had a handful of Word-proofing artefacts (<w:proofErr/>
# spell/grammar-check markers) left over from earlier edits, splitting
# merge-field placeholders such as <<courtName>> across multiple runs.
# This pass:
#   1) fixes the actual typo ("fee off" -> "fee of"), and
#   2) tidies the proofing markers around the placeholders, merging the
#      runs they used to straddle back into clean runs (mirroring what
#      Word itself does when the text is revisited/re-proofed).
#
# Helper: re-find a (unique) span of text, replacing it with itself (or
# with new text) via Find/Replace - this clears any <w:proofErr/> marks
# that lie fully inside the matched span and merges the runs the span
# touches into a single run.
function Merge-Range($doc, [string]$findText, [string]$replaceText) {
    $doc.Content.Find.Execute($findText, $true, $false, $false, $false, $false,
                               $true, 1, $false, $replaceText, 2) | Out-Null
}

# Helper: re-locate the (now merged) span and split it back into runs at
# the given internal offsets by toggling Bold on/off across each piece -
# a no-op to the visible formatting, but it forces Word to lay the piece
# down as its own <w:r>, with the run boundaries we want.
function Split-Range($doc, [string]$findText, [int[]]$lengths) {
    $rng = $doc.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
    $start = $rng.Start
    $pos = $start
    foreach ($len in $lengths) {
        $piece = $doc.Range($pos, $pos + $len)
        $piece.Font.Bold = 1
        $piece.Font.Bold = 0
        $pos = $pos + $len
    }
}

$d = $word.ActiveDocument

# --- Section 1: <<{dateFormat($nowUTC,'d MMMM yyyy')}>> -------------------
$full1 = "<<{dateFormat(`$nowUTC,`u2018d MMMM yyyy`u2019)}>>"
$full1 = "<<{dateFormat(`$nowUTC," + [char]0x2018 + "d MMMM yyyy" + [char]0x2019 + ")}>>"
Merge-Range $d $full1 $full1
Split-Range $d $full1 @(20, 1, 18)

# --- Section 2: <<courtName>> ---------------------------------------------
Merge-Range $d "<<courtName>>" "<<courtName>>"

# --- Section 3: <<claimNumber>> -------------------------------------------
Merge-Range $d "<<claimNumber>>" "<<claimNumber>>"
Split-Range $d "<<claimNumber>>" @(3, 4, 8)

# --- Section 4: <<applicationCreatedDate>> --------------------------------
Merge-Range $d "<<applicationCreatedDate>>" "<<applicationCreatedDate>>"
Split-Range $d "<<applicationCreatedDate>>" @(2, 22, 2)

# --- Section 5: typo fix + <<additionalApplicationFee>> -------------------
$feeFindOld = "be made with notice so that a hearing may be held. This means you will have to pay an additional fee off <<additionalApplicationFee>>."
$feeNew = "be made with notice so that a hearing may be held. This means you will have to pay an additional fee of <<additionalApplicationFee>>."
Merge-Range $d $feeFindOld $feeNew
Split-Range $d $feeNew @(116, 11, 6)
